$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ledger Generation Date: 2020-09-17 -> 2020-09-24 (keep as plain text, not an auto-converted date serial)
$ws.Range("B5").Value = "'2020-09-24"

# Row 10: 2020-09-10/test/Exp-1/EXPENSE/25 -> 2020-09-22/expense/Exp-1/EXPENSE/0
$ws.Range("A10").Value = "'2020-09-22"
$ws.Range("B10").Value = "expense"
$ws.Range("E10").Value = 0

# Row 11: 2020-09-10/test/Exp-2/EXPENSE/5 -> 2020-09-22/TEST/CN-0001/CREDITNOTE/0
$ws.Range("A11").Value = "'2020-09-22"
$ws.Range("B11").Value = "TEST"
$ws.Range("C11").Value = "CN-0001"
$ws.Range("D11").Value = "CREDITNOTE"
$ws.Range("E11").Value = 0

# The old TOTAL row (row 12) shifts down to row 13 - move its formatting + content first.
$ws.Range("D12:E12").Copy()
$ws.Range("D13:E13").PasteSpecial(-4122)
$ws.Range("D13").Value = "TOTAL"
$ws.Range("E13").Formula = "=SUM(E10:E12)"

# Clear the old TOTAL cells out of row 12 (they are replaced by new data below).
$ws.Range("D12:E12").ClearContents()
$ws.Range("D12:E12").ClearFormats()

# New row 12: second credit note, matching the style used for the data rows (10-11).
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Range("A12").Value = "'2020-09-23"
$ws.Range("B12").Value = "TEST"
$ws.Range("C12").Value = "CN-0002"
$ws.Range("D12").Value = "CREDITNOTE"
$ws.Range("E12").Value = 0
